# Update EPEX spot prices workbook with the latest day of data.
#
#   - "Prix Spot" sheet: append a new date column (DD) with 29-sep header
#     and 24 hourly prices, mirroring the style of the previous column.
#   - "Gaz" sheet: append two new daily rows (2025-09-27, 2025-09-28).
#   - "CO2" sheet: append two new daily rows (2025-09-27, 2025-09-28).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": add column DD (29-sep) after existing column DC.
# ---------------------------------------------------------------------
$spot = $wb.Worksheets.Item("Prix Spot")

# Clone the header cell's formatting (bold, centered, boxed) from DC1
# onto DD1 without disturbing DD1's eventual value.
$spot.Range("DC1").Copy() | Out-Null
$spot.Range("DD1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$spot.Range("DD1").Value = "29-sep"

$hourlyPrices = @(
    51.6,
    47.7,
    45,
    32,
    30,
    32.5,
    50,
    74.40000000000001,
    81.8,
    73.08,
    59.33,
    40.14,
    35,
    28.3,
    21.99,
    27.94,
    35,
    52,
    85,
    101.13,
    133.26,
    89.51000000000001,
    84.88,
    88.59999999999999
)

for ($i = 0; $i -lt $hourlyPrices.Count; $i++) {
    $row = $i + 2
    $spot.Cells.Item($row, 108).Value = $hourlyPrices[$i]
}

# ---------------------------------------------------------------------
# Sheet "Gaz": append 2025-09-27 and 2025-09-28 rows.
# ---------------------------------------------------------------------
$gaz = $wb.Worksheets.Item("Gaz")
$gaz.Range("A105").Value = "'2025-09-27"
$gaz.Range("B105").Value = 31.775
$gaz.Range("A106").Value = "'2025-09-28"
$gaz.Range("B106").Value = 31.775

# ---------------------------------------------------------------------
# Sheet "CO2": append 2025-09-27 and 2025-09-28 rows.
# ---------------------------------------------------------------------
$co2 = $wb.Worksheets.Item("CO2")
$co2.Range("A105").Value = "'2025-09-27"
$co2.Range("B105").Value = 75.26000000000001
$co2.Range("A106").Value = "'2025-09-28"
$co2.Range("B106").Value = 75.26000000000001
